# Apply the edits described by the diff:
#  - Update "untere Grenze [MW]" (col D) values from 99 to 100 for several rows
#  - Update "obere Grenze [MW]" (col E) values for rows 14 and 15
#  - Update the active selection on the sheet to N13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Column D ("untere Grenze [MW]") updates: 99 -> 100
$ws.Range("D12").Value = 100
$ws.Range("D13").Value = 100
$ws.Range("D16").Value = 100
$ws.Range("D17").Value = 100
$ws.Range("D18").Value = 100
$ws.Range("D19").Value = 100
$ws.Range("D20").Value = 100
$ws.Range("D21").Value = 100
$ws.Range("D22").Value = 100
$ws.Range("D23").Value = 100
$ws.Range("D24").Value = 100
$ws.Range("D25").Value = 100

# Column E ("obere Grenze [MW]") updates
$ws.Range("E14").Value = 99999
$ws.Range("E15").Value = 9999

# Update the selected cell/range on the active sheet
$ws.Activate()
$ws.Range("N13").Select()
